$wb = $excel.ActiveWorkbook

# --- Sheet "Leads": insert a new lead (Martha Bixby) above the existing
# Marty McFly row, which shifts down to row 3 (with its hyperlink style).
$leads = $wb.Worksheets.Item("Leads")

$leads.Rows.Item(2).Insert()

# The row insert carries the hyperlink *style* down to D3, but this
# runtime leaves the hyperlink collection entry pointing at the old (now
# blank) D2 cell, so drop it before re-adding it against D3.
$leads.Range("D2").Hyperlinks.Delete()

# Fill in the new lead's data (order matters for shared-string layout).
$leads.Range("A2").Value = "Martha"
$leads.Range("B2").Value = "Bixby"
$leads.Range("D2").Value = "martha.bixby@demo.mail"
$leads.Range("E2").Value = "Web"
$leads.Range("F2").Value = "Private"
$leads.Range("C2").Value = "1-222-505-424"

# Re-create the hyperlink for Marty's e-mail address on its new row, and
# give the e-mail-style "Hyperlink" cell format to both e-mail cells
# (Martha's keeps the visual style without being a clickable link).
$leads.Hyperlinks.Add($leads.Range("D3"), "mailto:mctester@deloreantesting.com") | Out-Null
$leads.Range("D3").Style = "Hyperlink"
$leads.Range("D2").Style = "Hyperlink"

$leads.Columns.Item(3).ColumnWidth = 12.53125

$leads.Range("C2").Select() | Out-Null

# --- Sheet "Cars": the LeadName column is no longer used.
$cars = $wb.Worksheets.Item("Cars")
$cars.Columns.Item(4).ClearContents()

$cars.Range("D1:D1048576").Select() | Out-Null

$leads.Activate() | Out-Null
